# Actualización automática 2025-08-13 11:15:10
# Registers a sale of 142.56 (PIEDRA SINTERIZADA, agosto) for client
# "ARCOS GOMEZ CONSTRUCCIONES CIA. LTDA." and refreshes the dependent
# summary cells across the three sheets of the workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO": per-group sales for ARCOS GOMEZ (row 3) ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("L3").Value = 142.56          # PIEDRA SINTERIZADA
$wsGrupo.Range("L12").Value = "1 de 10"      # clients-with-sales counter

# --- Sheet "VENTA MENSUAL": monthly sales for ARCOS GOMEZ (row 3) ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F3").Value = 142.56        # agosto
$wsMensual.Range("F12").Value = 142.56       # agosto column total

# --- Sheet "CUMPLIMIENTO MENSUAL": compliance roll-up ---
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
# Excel's ColumnWidth (chars) <-> stored OOXML width has a ~0.83-char
# offset baked in by the Normal-style font metrics; 11.17 is what round-trips
# to the target stored width of 12.
$wsCumpl.Columns.Item(4).ColumnWidth = 11.17

$wsCumpl.Range("D2").Value = 142.56
$wsCumpl.Range("E2").Value = -142.56

$wsCumpl.Range("D4").Value = 142.56
$wsCumpl.Range("E4").Value = 9857.440000000001
$wsCumpl.Range("F4").Value = 0.014256
